$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the abbreviated party acronyms in the header row (B1:O1) into
# their full descriptive names.
$ws.Range("B1").Value = "CD - The Centre Democrats (Centrum-Demokraterne , CD)"
$ws.Range("C1").Value = "EL - The Unity List (Enhedslisten , EL)"
$ws.Range("D1").Value = "FRP - The Progress Party (Fremskridtspartiet , FRP)"
$ws.Range("E1").Value = "KF - The Conservative People’s Party ( Det Konservative Folkeparti , KF)"
$ws.Range("F1").Value = "RV - The Social-Liberal Party (Det radikale Venstre , RV)"
$ws.Range("G1").Value = "SD - The Social Democratic Party (Socialdemokratiet , SD)"
$ws.Range("H1").Value = "SF - The Socialist People’s Party  (Socialistisk Folkeparti , SF)"
$ws.Range("I1").Value = "V - The Liberal Party (Venstre , V)"
$ws.Range("J1").Value = "no acronym - Independents (Uafhængig, no acronym)"
$ws.Range("K1").Value = "DF - Danish People's Party  (Dansk Folkeparti , DF)"
$ws.Range("L1").Value = "KD - Christian Democrats  (Kristendemokraterne , KD), known until 2004 as The Christian People’s Party  (KRF, Kristeligt Folkeparti )"
$ws.Range("M1").Value = "LA - The Liberal Alliance   (The Liberal Alliance , LA), known until 27 August 2008 as New Alliance (NA, Ny Alliance )"
$ws.Range("N1").Value = "Å - The Alternative (Alternativet, Å)"
$ws.Range("O1").Value = "NB - New Right (Nye Borgerlige, NB)"

# Clean up floating-point rounding artifacts in the computed minimal-seat
# values so they store as exact integers.
$ws.Range("D3").Value = 2
$ws.Range("F3").Value = 3
$ws.Range("K3").Value = 7
$ws.Range("E9").Value = 5
$ws.Range("H9").Value = 6
$ws.Range("M9").Value = 2
$ws.Range("N9").Value = 2
$ws.Range("O9").Value = 2
